$wb = $excel.ActiveWorkbook

# --- county-year sheet: row for year 2002 (row 6) ---
# incorporate new urban returns, use 2000 data for 2002
$wsCountyYear = $wb.Worksheets.Item("county-year")
$wsCountyYear.Range("B6").Value = 2069
$wsCountyYear.Range("C6").Value = 67.350257873535156
$wsCountyYear.Range("S6").Value = 3072
$wsCountyYear.Range("T6").Value = 100
$wsCountyYear.Range("U6").Value = 0
$wsCountyYear.Range("V6").Value = 0

# --- county-overall sheet: overall totals row (row 2) ---
$wsCountyOverall = $wb.Worksheets.Item("county-overall")
$wsCountyOverall.Range("A2").Value = 6123
$wsCountyOverall.Range("B2").Value = 24.91455078125
$wsCountyOverall.Range("R2").Value = 12227
$wsCountyOverall.Range("S2").Value = 49.751789093017578
$wsCountyOverall.Range("T2").Value = 12349
$wsCountyOverall.Range("U2").Value = 50.248210906982422

# --- point-year sheet: row for year 2002 (row 6) ---
$wsPointYear = $wb.Worksheets.Item("point-year")
$wsPointYear.Range("B6").Value = 929707
$wsPointYear.Range("C6").Value = 68.2293701171875
$wsPointYear.Range("S6").Value = 1362620
$wsPointYear.Range("T6").Value = 100
$wsPointYear.Range("U6").Value = 0
$wsPointYear.Range("V6").Value = 0

# --- point-overall sheet: overall totals row (row 2) ---
$wsPointOverall = $wb.Worksheets.Item("point-overall")
$wsPointOverall.Range("A2").Value = 2749941
$wsPointOverall.Range("B2").Value = 25.226594924926758
$wsPointOverall.Range("R2").Value = 5417373
$wsPointOverall.Range("S2").Value = 49.696292877197266
$wsPointOverall.Range("T2").Value = 5483587
$wsPointOverall.Range("U2").Value = 50.303707122802734
